$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list on Mon Apr  1 07:28:03 UTC 2024 with GitHub Actions
#
# Each target cell is temporarily switched to Text number format before
# its value is written, so that numeric-looking strings (e.g. "2.80",
# "69.646.36", "0.630") retain their exact textual representation instead
# of being coerced into floating point numbers by Excel. The style is then
# reset back to "Normal" so the cell formatting matches the original
# (unstyled) cells.

$ws.Range("D2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.646.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.546.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.83%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "586.62"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.05%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.05%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.630"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.41%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -5.09%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.34%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "687.83"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +16.16%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.110.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.741.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.554.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.49%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.972"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.92"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.58%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "108.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +5.46%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.74%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.79%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.83%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.69%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.05%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.68"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.93"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.98%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.93"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -4.10%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "62.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.799.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.54%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -8.03%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +4.24%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "501.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.82%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.35%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.136"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.56%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.22%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "34.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.75%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0459"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.65%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.62%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Style = "Normal"
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Style = "Normal"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +71.12%  "
$ws.Range("E51").Style = "Normal"
